# Aggiornamento fino al 9 agosto 2021
# Append rows 329-343 (dates 44403-44417) to Sheet1, mirroring the
# formatting of the preceding rows (column A carries the date style,
# columns B/C/D remain unstyled numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 329

# date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$data = @(
    @(44403, 0, 3, 49.90850108135086),
    @(44404, 0, 3, 49.90850108135086),
    @(44405, 0, 3, 49.90850108135086),
    @(44406, 1, 3, 49.90850108135086),
    @(44407, 0, 3, 49.90850108135086),
    @(44408, 0, 1, 16.63616702711695),
    @(44409, 0, 1, 16.63616702711695),
    @(44410, 0, 1, 16.63616702711695),
    @(44411, 0, 1, 16.63616702711695),
    @(44412, 0, 1, 16.63616702711695),
    @(44413, 1, 1, 16.63616702711695),
    @(44414, 0, 1, 16.63616702711695),
    @(44415, 0, 1, 16.63616702711695),
    @(44416, 0, 1, 16.63616702711695),
    @(44417, 0, 1, 16.63616702711695)
)

$endRow = $startRow + $data.Length - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

# Copy the formatting (style/borders/bold date format) from the last
# existing data row (328) onto the newly added rows so the new cells
# match the established look of the sheet.
$ws.Range("A328").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)

$ws.Range("B328:D328").Copy()
$ws.Range("B$startRow`:D$endRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0
